$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All of the "content" cells in columns B/C end up holding text that already
# exists elsewhere in the sheet, so each change is done as a cell-range copy
# (this keeps the shared-string table / cell styles byte-identical to a plain
# value assignment, and also avoids Excel's automatic text->date conversion
# for the "01/01/2022" value). Copies that read a row which is itself a copy
# target further down run first (bottom of the dependency chain up), so every
# source is read before it gets overwritten.

$ws.Range("B23:C23").Copy($ws.Range("B22:C22"))
$ws.Range("B22:C22").Copy($ws.Range("B21:C21"))
$ws.Range("B21:C21").Copy($ws.Range("B20:C20"))
$ws.Range("B20:C20").Copy($ws.Range("B14:C14"))

$ws.Range("B13:C13").Copy($ws.Range("B10:C10"))
$ws.Range("B13:C13").Copy($ws.Range("B17:C17"))
$ws.Range("B8:C8").Copy($ws.Range("B15:C15"))

# The two "docente" rows (13 and 14) held only B/C values (no label in column
# A); removing them shifts every row below up by two, which realigns all of
# the remaining A-column labels with their correct (new) row numbers and
# keeps the per-row heights in sync automatically.
$ws.Rows("13:14").Delete()
